# Adds the new localization rows (Video/Audio/Controls settings screens,
# pause-menu strings, volume-mixer labels, etc.) to the end of the
# en / zh localization table on Sheet1.
#
# Columns: A = key (en, used as lookup key), B = en display text,
#          C = zh display text.  Existing rows occupy 1:19; new rows
#          are appended starting at row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Video",            "Video",            "图像"),
    @("Audio",             "Audio",            "音频"),
    @("Controls",          "Controls",         "操作"),
    @("Return",            "Return",           "返回"),
    @("Fullscreen",        "Fullscreen",       "全屏"),
    @("Borderless",        "Borderless",       "无边框"),
    @("Disabled",          "Disabled",         "关闭"),
    @("Enabled",           "Enabled",          "开启"),
    @("Adaptive",          "Adaptive",         "自适应"),
    @("Nvidia DLSS",       "Nvidia DLSS",      "Nvidia DLSS"),
    @("Master_Music",      "Master",           "主音量"),
    @("Music",             "Music",            "音乐"),
    @("SFX",               "SFX",              "音效"),
    @("Game",              "Game",             "游戏"),
    @("Paused_Settings",   "settings",         "设置"),
    @("Resume",            "Resume",           "恢复游戏"),
    @("Restart",           "Restart",          "重新开始"),
    @("Main Menu",         "Main Menu",        "主菜单"),
    @("Game Paused",       "Game Paused",      "游戏暂停")
)

$startRow = 20
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $item = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
}

$lastRow = $startRow + $newRows.Count - 1

# Mirror the author's final selection/scroll state (best effort - window
# chrome position is host-managed and may not round-trip through OOXML).
$ws.Range("C$lastRow").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
